$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped entirely from the data set
# (row 28 = "SC 92", row 26 = "RM 232"). Delete from the bottom up so
# the remaining row indices used below are not disturbed mid-way.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# After the two rows above are removed, apply the individual cell-level
# value changes (new values, or clearing existing ones) to reach the
# final data state.
$ws.Range("E2").ClearContents()
$ws.Range("F4").ClearContents()

$ws.Range("D6").Value = -14.2
$ws.Range("F6").Value = 16.43

$ws.Range("D8").ClearContents()

$ws.Range("F12").ClearContents()

$ws.Range("F14").Value = 17.76

$ws.Range("D18").Value = -15.2

$ws.Range("D20").ClearContents()

$ws.Range("F21").Value = 16.58

$ws.Range("F22").Value = 16.81

$ws.Range("D23").Value = -13.9

$ws.Range("D25").ClearContents()

$ws.Range("F26").ClearContents()

$ws.Range("C27").Value = 10
$ws.Range("F27").ClearContents()

$ws.Range("C28").ClearContents()
$ws.Range("F28").ClearContents()

$ws.Range("C29").ClearContents()

$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7

$ws.Range("F31").Value = 17.18

$ws.Range("C32").ClearContents()
